$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data (2021-10-07) needs to be inserted right after the
# existing row 119, pushing every row from the old row 120 onward down by two
# rows. Insert two blank rows at row 120 to achieve that shift.
$ws.Rows(120).Insert()
$ws.Rows(120).Insert()

# Populate the two new rows with the new week's data. Columns A, B, C, E, F,
# G, I, N, O, Q, R repeat the same boilerplate values used throughout the
# sheet for this market/product.

# Row 120 - Zafiro rojo
$ws.Range("A120").Value = 7
$ws.Range("B120").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C120").Value = "Ñuble"
$ws.Range("D120").Value = 44476
$ws.Range("E120").Value = 16
$ws.Range("F120").Value = 100112002
$ws.Range("G120").Value = "Pimiento"
$ws.Range("H120").Value = "Zafiro rojo"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 160
$ws.Range("K120").Value = 43000
$ws.Range("L120").Value = 44000
$ws.Range("M120").Value = 43500
$ws.Range("N120").Value = "$/caja 15 kilos"
$ws.Range("O120").Value = "Región de Arica y Parinacota"
$ws.Range("P120").Value = 2900
$ws.Range("Q120").Value = 15
$ws.Range("R120").Value = "Hortaliza"

# Row 121 - Zafiro verde
$ws.Range("A121").Value = 7
$ws.Range("B121").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C121").Value = "Ñuble"
$ws.Range("D121").Value = 44476
$ws.Range("E121").Value = 16
$ws.Range("F121").Value = 100112002
$ws.Range("G121").Value = "Pimiento"
$ws.Range("H121").Value = "Zafiro verde"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 160
$ws.Range("K121").Value = 41000
$ws.Range("L121").Value = 42000
$ws.Range("M121").Value = 41500
$ws.Range("N121").Value = "$/caja 15 kilos"
$ws.Range("O121").Value = "Región de Arica y Parinacota"
$ws.Range("P121").Value = 2767
$ws.Range("Q121").Value = 15
$ws.Range("R121").Value = "Hortaliza"
